$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 14552
$ws1.Range("F4").Value = 699
$ws1.Range("F6").Value = 599
$ws1.Range("F7").Value = 1543

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 14552
$ws4.Range("F4").Value = 699
$ws4.Range("F8").Value = 599
$ws4.Range("F9").Value = 1543
